# "Add files via upload" - estadisticas_2026.xlsx
# Update Jugadores (player position correction), add the newest Partido
# (match 5) to Partidos, and add its Eventos (event log) rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Jugadores: player id 19 (JUAN DIEGO GOMEZ CEBALLOS) is now listed
#    as "mediocampista" instead of "defensa".
# ---------------------------------------------------------------------
$wsJugadores = $wb.Worksheets.Item("Jugadores")
$wsJugadores.Range("C40").Value = "mediocampista"

# ---------------------------------------------------------------------
# 2) Partidos: append the 5th match played on 2026-01-31 at "El Seminario".
# ---------------------------------------------------------------------
$wsPartidos = $wb.Worksheets.Item("Partidos")

$wsPartidos.Cells.Item(6, 1).Value = 5
$wsPartidos.Cells.Item(6, 2).Value = 46053
$wsPartidos.Cells.Item(6, 3).Value = "g"
$wsPartidos.Cells.Item(6, 4).Value = "p"
$wsPartidos.Cells.Item(6, 5).Value = 7
$wsPartidos.Cells.Item(6, 6).Value = 5
$wsPartidos.Cells.Item(6, 7).Value = "El Seminario"

# Match date column uses a date number format - copy it from the row above
# instead of the default General format.
$wsPartidos.Range("B5").Copy()
$wsPartidos.Range("B6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Eventos: the two previous "azul"-team placeholder rows (72-73) no
#    longer carry the leftover border formatting.
# ---------------------------------------------------------------------
$wsEventos = $wb.Worksheets.Item("Eventos")
$wsEventos.Range("A72:I73").Style = "Normal"

# Append the 22 player-event rows belonging to match 5 (rows 102-123).
$eventRows = @(
    @(5,51,"azul",7,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,88,"azul",0,0,0,0,0,1,0,0,0,1,0,0,1),
    @(5,95,"azul",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,15,"azul",0,0,0,0,0,0,0,0,0,1,0,0,1),
    @(5,29,"azul",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,7,"azul",0,0,0,0,0,0,1,0,1,0,0,0,1),
    @(5,71,"azul",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,21,"azul",0,0,0,1,0,0,0,0,1,0,0,0,1),
    @(5,10,"azul",0,0,0,0,0,0,0,0,1,0,0,0,1),
    @(5,6,"azul",0,0,0,0,0,1,1,0,0,0,0,0,1),
    @(5,14,"azul",0,0,0,0,0,0,1,0,1,0,0,0,1),
    @(5,1,"amarillo",5,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,3,"amarillo",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,26,"amarillo",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,19,"amarillo",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,27,"amarillo",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,80,"amarillo",0,0,0,0,0,1,0,0,1,0,0,0,1),
    @(5,52,"amarillo",0,0,0,0,0,0,1,0,0,0,0,0,1),
    @(5,24,"amarillo",0,0,0,0,0,0,0,0,0,0,0,0,1),
    @(5,73,"amarillo",0,0,0,0,0,2,2,0,0,0,0,0,1),
    @(5,12,"amarillo",0,0,0,0,0,0,1,0,2,0,0,0,1),
    @(5,13,"amarillo",0,0,0,0,0,1,0,0,1,0,0,0,1)
)

$startRow = 102
for ($i = 0; $i -lt $eventRows.Count; $i++) {
    $r = $startRow + $i
    $row = $eventRows[$i]

    $wsEventos.Cells.Item($r, 1).Value = $row[0]   # id_partido
    $wsEventos.Cells.Item($r, 2).Value = $row[1]   # id_jugador
    $wsEventos.Cells.Item($r, 3).Value = $row[2]   # equipo
    $wsEventos.Cells.Item($r, 4).Value = $row[3]   # gol_recibido
    $wsEventos.Cells.Item($r, 5).Value = $row[4]   # fue_arquero
    $wsEventos.Cells.Item($r, 6).Value = $row[5]   # fue_defensa
    $wsEventos.Cells.Item($r, 7).Value = $row[6]   # fue_mediocampista
    $wsEventos.Cells.Item($r, 8).Value = $row[7]   # fue_delantero
    $wsEventos.Cells.Item($r, 9).Value = $row[8]   # gol_primer
    $wsEventos.Cells.Item($r, 10).Value = $row[9]  # gol_segundo
    $wsEventos.Cells.Item($r, 11).Formula = "=I" + $r + "+J" + $r   # gol_total
    $wsEventos.Cells.Item($r, 12).Value = $row[10] # autogoles
    $wsEventos.Cells.Item($r, 13).Value = $row[11] # asistencia_gol
    $wsEventos.Cells.Item($r, 14).Value = $row[12] # amarillas
    $wsEventos.Cells.Item($r, 15).Value = $row[13] # rojas
    $wsEventos.Cells.Item($r, 16).Value = $row[14] # penal_atajado
    $wsEventos.Cells.Item($r, 17).Value = $row[15] # partido_completado
}

# ---------------------------------------------------------------------
# 4) Restore each sheet's on-screen selection to match where the author
#    left off while editing. "Eventos" is re-activated last so it stays
#    the workbook's active tab, as it was before the edit.
# ---------------------------------------------------------------------
$wsJugadores.Activate()
$wsJugadores.Range("B28").Select()

$wsPartidos.Activate()
$wsPartidos.Range("G7").Select()

$wsEventos.Activate()
$wsEventos.Range("Q113").Select()
